$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "1" to "ადიგენი"
$ws.Name = "ადიგენი"

# Clear the census-note text (and its cell formatting) in A2, leaving a blank spacer row
$ws.Range("A2").Clear()

# Remove the now-empty spacer row (old row 3) so everything below shifts up
$ws.Rows("3:3").Delete()

# Remove the 1989 and 2002 data columns (old B and C), leaving only the 2014 column
$ws.Columns("B:C").Delete()

# The 2014 header cell no longer needs the heavier outer-edge border that used to mark
# the right-hand boundary of the whole table; give it the plain thin divider instead
$ws.Range("B4").Borders.Item([Microsoft.Office.Interop.Excel.XlBordersIndex]::xlEdgeRight).Weight = [Microsoft.Office.Interop.Excel.XlBorderWeight]::xlThin

# Select A2 to match the recorded active selection in the target file
$ws.Range("A2").Select()
